# Apply the crypto price/volume updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.409.53"
$ws.Range("E2").Value = "  -5.76%  "

$ws.Range("D3").Value = "3.460.87"
$ws.Range("E3").Value = "  -7.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "557.01"
$ws.Range("E5").Value = "  -9.08%  "

$ws.Range("D6").Value = "180.18"
$ws.Range("E6").Value = "  -6.49%  "

$ws.Range("D7").Value = "3.453.99"

$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -6.66%  "

$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "0.642"
$ws.Range("E10").Value = "  -12.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.140"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -13.93%  "

$ws.Range("D12").Value = "50.86"
$ws.Range("E12").Value = "  -16.20%  "

$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  -14.81%  "

$ws.Range("D14").Value = "9.39"
$ws.Range("E14").Value = "  -12.27%  "

$ws.Range("D15").Value = "4.006.73"
$ws.Range("E15").Value = "  -7.33%  "

$ws.Range("E16").Value = "  -1.90%  "

$ws.Range("D17").Value = "3.450.44"
$ws.Range("E17").Value = "  -7.56%  "

$ws.Range("D18").Value = "65.095.94"
$ws.Range("E18").Value = "  -5.93%  "

$ws.Range("D19").Value = "17.58"
$ws.Range("E19").Value = "  -10.10%  "

$ws.Range("D20").Value = "11.57"
$ws.Range("E20").Value = "  -10.93%  "

$ws.Range("E21").Value = "  -11.65%  "

$ws.Range("D22").Value = "374.98"
$ws.Range("E22").Value = "  -9.55%  "

$ws.Range("D23").Value = "4.05"
$ws.Range("E23").Value = "  -11.74%  "

$ws.Range("D24").Value = "82.15"
$ws.Range("E24").Value = "  -8.68%  "

$ws.Range("D25").Value = "10.59"
$ws.Range("E25").Value = "  -3.70%  "

$ws.Range("E26").Value = "  -1.18%  "

$ws.Range("D27").Value = "2.77"
$ws.Range("E27").Value = "  -10.05%  "

$ws.Range("D28").Value = "11.75"
$ws.Range("E28").Value = "  -9.24%  "

$ws.Range("D29").Value = "3.42"
$ws.Range("E29").Value = "  -10.89%  "

$ws.Range("D30").Value = "8.54"
$ws.Range("E30").Value = "  -12.49%  "

$ws.Range("D31").Value = "29.98"
$ws.Range("E31").Value = "  -9.85%  "

$ws.Range("D32").Value = "7.05"
$ws.Range("E32").Value = "  -9.48%  "

$ws.Range("D33").Value = "607.25"
$ws.Range("E33").Value = "  -4.67%  "

$ws.Range("D34").Value = "11.74"
$ws.Range("E34").Value = "  -8.45%  "

$ws.Range("D35").Value = "62.42"
$ws.Range("E35").Value = "  -5.54%  "

$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  -12.72%  "

$ws.Range("D37").Value = "40.18"
$ws.Range("E37").Value = "  -12.54%  "

$ws.Range("E38").Value = "  +0.07%  "

$ws.Range("D39").Value = "0.392"
$ws.Range("E39").Value = "  -6.29%  "

$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("D41").Value = "0.0₃0705"
$ws.Range("E41").Value = "  -15.98%  "

$ws.Range("E42").Value = "  -9.52%  "

$ws.Range("D43").Value = "2.883.37"
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("D44").Value = "2.69"
$ws.Range("E44").Value = "  -12.76%  "

$ws.Range("E45").Value = "  -8.62%  "

$ws.Range("D46").Value = "3.07"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0390"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -13.45%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.126"
$ws.Range("E48").Value = "  -10.41%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "137.05"
$ws.Range("E49").Value = "  -3.81%  "

$ws.Range("E50").Value = "  -11.61%  "

$ws.Range("D51").Value = "8.04"
$ws.Range("E51").Value = "  -12.59%  "
